# Trade #50 closed at 2026-02-17 08:39:39 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up figures for the
# MarketMaking strategy after trade #50 closed, and appends the new
# trade row to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet - refreshed account-level totals
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.51   # Current Capital
$summary.Range("B4").Value = -0.49     # Total P&L $
$summary.Range("B5").Value = -0.2      # Total P&L %
$summary.Range("B6").Value = 50        # Total Trades
$summary.Range("B8").Value = 22        # Losing Trades
$summary.Range("B9").Value = 36        # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.51000000000001   # Capital
$status.Range("D4").Value = 50                  # Trades
$status.Range("E4").Value = -0.49               # P&L $
$status.Range("F4").Value = -0.49               # P&L %
$status.Range("G4").Value = 36                  # Win Rate %

# ---------------------------------------------------------------
# 3. Append the closed trade as row 51 to both trade logs
# ---------------------------------------------------------------
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 51

    $ws.Cells.Item($row, 1).Value = 50                  # A - Trade #

    # Date/Time columns must stay plain text (like the rest of the
    # log) instead of being auto-converted to Excel date/time
    # serials, so force a text format before writing them.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"        # B - Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:39:32"          # C - Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"      # D - Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"              # E - Side
    $ws.Cells.Item($row, 6).Value = 0.14                # F - Entry Price
    $ws.Cells.Item($row, 7).Value = 0.1                 # G - Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"            # H - Status
    $ws.Cells.Item($row, 9).Value = -28.5714            # I - P&L %
    $ws.Cells.Item($row, 10).Value = -0.04              # J - P&L $
    $ws.Cells.Item($row, 11).Value = 99.51000000000001  # K - Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # L - Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # M - Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # N - Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O - Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"       # P - Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13                # Q - Duration (min)
}

Write-Output "Applied trade #50 close updates"
